$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 42 (add missing tested/hosp/icu/death figures) ---
$ws.Range("B42").Value = 1248
$ws.Range("AA42").Value = 480
$ws.Range("AB42").Value = 170
$ws.Range("AC42").Value = 60

# --- New row 43 (2020-04-16, serial 43937) ---
$ws.Range("A43").Value = 43937
$ws.Range("B43").Value = 1114
$ws.Range("C43").Value = 2087
$ws.Range("D43").Value = 12
$ws.Range("F43").Value = 34
$ws.Range("H43").Value = 316
$ws.Range("J43").Value = 399
$ws.Range("L43").Value = 362
$ws.Range("N43").Value = 383
$ws.Range("P43").Value = 286
$ws.Range("R43").Value = 163
$ws.Range("T43").Value = 129
$ws.Range("V43").Value = 3
$ws.Range("X43").Value = 1034
$ws.Range("Y43").Value = 1047
$ws.Range("Z43").Value = 6
$ws.Range("AA43").Value = 507
$ws.Range("AB43").Value = 181
$ws.Range("AC43").Value = 63

# --- New row 44 (2020-04-17, serial 43938) ---
$ws.Range("A44").Value = 43938
$ws.Range("C44").Value = 2158
$ws.Range("D44").Value = 13
$ws.Range("F44").Value = 35
$ws.Range("H44").Value = 327
$ws.Range("J44").Value = 408
$ws.Range("L44").Value = 373
$ws.Range("N44").Value = 404
$ws.Range("P44").Value = 297
$ws.Range("R44").Value = 166
$ws.Range("T44").Value = 131
$ws.Range("V44").Value = 4
$ws.Range("X44").Value = 1073
$ws.Range("Y44").Value = 1081
$ws.Range("Z44").Value = 4
$ws.Range("AA44").Value = 524
$ws.Range("AB44").Value = 181
$ws.Range("AC44").Value = 70

# Move the active selection to the last-entered cell, matching the
# author's final cursor position after typing in the new data.
[void]$ws.Range("B44").Select()
